$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Module parameters update: DTC (row 33), FD (row 34), A (row 35), B (row 36) ---

# DTC: "Sandia" (col B) and "mono-Si" (col D) values updated to 3; right-align col D to match B/C
$ws.Range("B33").Value = 3
$ws.Range("D33").Value = 3
$ws.Range("D33").HorizontalAlignment = -4152  ## xlRight

# FD: only alignment of col D needs to match B/C now (value unchanged, still 1)
$ws.Range("D34").HorizontalAlignment = -4152  ## xlRight

# A: col B and col D updated to -3.56
$ws.Range("B35").Value = -3.56
$ws.Range("D35").Value = -3.56
$ws.Range("D35").HorizontalAlignment = -4152  ## xlRight

# B: col B and col D updated to -0.075
$ws.Range("B36").Value = -0.074999999999999997
$ws.Range("D36").Value = -0.074999999999999997
$ws.Range("D36").HorizontalAlignment = -4152  ## xlRight

# --- View state: scroll to column C and select G11 ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G11").Select()
